$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# --- Text / value updates -------------------------------------------------
# "In Translation" -> "Ready for handoff" (Status columns)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Updated generation / handoff timestamps
$wsOverview.Range("G2").Value = "2016-08-20 17:06:15"
$wsDeDe.Range("H2").Value = "2016-08-20 17:06:15"
$wsZhCn.Range("H2").Value = "2016-08-20 17:06:11"

# --- Column width updates ---------------------------------------------------
# Target stored width is 17.2159881591797; the COM layer quantizes
# ColumnWidth to 1/6-character increments, so feed it the value whose
# quantized result lands on the closest achievable grid point
# (stored width 17.166666666666668).
$newWidth = 16.38265482584637
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
